# Extend the "Your Kind Co-operation is Solicited to below mentioned
# account." sentence with ". Kindly issue and hand over the statement to
# Mr Krishna kumar Gupta.", splitting the addition into a new paragraph
# right after "Kindly issue and ".

$d = $word.ActiveDocument

# Locate the paragraph that holds the sentence we need to extend.
$targetPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -match "to below mentioned account") {
        $targetPara = $p
        break
    }
}

# Anchor on " to below mentioned account" inside that paragraph and
# collapse the range to right after it (i.e. just before the existing
# closing ".").
$r = $targetPara.Range
$r.Find.Execute("to below mentioned account", $true, $false, $false, $false, $false, `
                 $true, 1, $false, "", 0) | Out-Null
$r.Collapse(0)

# Insert the clause that stays in the first paragraph, right before the
# original trailing ".".
$r.InsertAfter(". Kindly issue and ")
$r.Collapse(0)

# Break the paragraph here; the original trailing "." run becomes the
# tail of the brand-new paragraph and automatically inherits the
# sz/szCs/lang run formatting via the paragraph mark.
$r.InsertParagraphAfter()

# Insert the new sentence text at the start of that new paragraph, i.e.
# immediately before the carried-over "." run.
$newParaStart = $r.End + 1
$r2 = $d.Range($newParaStart, $newParaStart)
$r2.InsertBefore("hand over the statement to Mr Krishna kumar Gupta")
